$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 58826108
$ws.Range("I28").Value = 76925970
$ws.Range("K28").Value = 76925970
$ws.Range("M28").Value = -76925485
$ws.Range("H43").Value = 1691
$ws.Range("I43").Value = 1691
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1691
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1622
$ws.Range("H51").Value = 3310.423
$ws.Range("J51").Value = 3185.0454
$ws.Range("L51").Value = 3185.0454
$ws.Range("N51").Value = -4153.0454
$ws.Range("H53").Value = 273.5909
$ws.Range("I53").Value = 320
$ws.Range("J53").Value = 217.9
$ws.Range("K53").Value = 320
$ws.Range("L53").Value = 217.9
$ws.Range("M53").Value = 317
$ws.Range("N53").Value = -1491.9
$ws.Range("H62").Value = 4833
$ws.Range("I62").Value = 4833
$ws.Range("K62").Value = 4833
$ws.Range("M62").Value = -4209
$ws.Range("H65").Value = 4833
$ws.Range("I65").Value = 4833
$ws.Range("K65").Value = 24165
$ws.Range("M65").Value = -21045
$ws.Range("H98").Value = 1286.9565
$ws.Range("I98").Value = 1247.5
$ws.Range("K98").Value = 1247.5
$ws.Range("M98").Value = 250.5
$ws.Range("H106").Value = 45460456
$ws.Range("I106").Value = 55560556
$ws.Range("K106").Value = 55560556
$ws.Range("M106").Value = -55559925
$ws.Range("H122").Value = 1286.9565
$ws.Range("I122").Value = 1247.5
$ws.Range("K122").Value = 3742.5
$ws.Range("M122").Value = -1292.5
$ws.Range("H129").Value = 1317.8572
$ws.Range("I129").Value = 787.75
$ws.Range("K129").Value = 2363.25
$ws.Range("M129").Value = 2636.75
$ws.Range("N43").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4418.213
$ws.Range("I32").Value = 2752.611
$ws.Range("J32").Value = 17267.143
$ws.Range("K32").Value = 2752.611
$ws.Range("L32").Value = 17267.143
$ws.Range("M32").Value = -2465.611
$ws.Range("N32").Value = -17841.143
$ws.Range("H45").Value = 13878.286
$ws.Range("I45").Value = 19762.25
$ws.Range("J45").Value = 6033
$ws.Range("K45").Value = 19762.25
$ws.Range("L45").Value = 6033
$ws.Range("M45").Value = -19385.25
$ws.Range("N45").Value = -6787
$ws.Range("H74").Value = 1475.375
$ws.Range("I74").Value = 1257.5714
$ws.Range("K74").Value = 1257.5714
$ws.Range("M74").Value = -383.5714
$ws.Range("H77").Value = 1475.375
$ws.Range("I77").Value = 1257.5714
$ws.Range("K77").Value = 6287.857
$ws.Range("M77").Value = -1919.857

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4349.375
$ws.Range("I16").Value = 3874.25
$ws.Range("J16").Value = 4824.5
$ws.Range("K16").Value = 3874.25
$ws.Range("L16").Value = 4824.5
$ws.Range("M16").Value = -3587.25
$ws.Range("N16").Value = -5398.5
$ws.Range("H19").Value = 2452.7058
$ws.Range("I19").Value = 1306
$ws.Range("J19").Value = 3255.4
$ws.Range("K19").Value = 1306
$ws.Range("L19").Value = 3255.4
$ws.Range("M19").Value = -1136
$ws.Range("N19").Value = -3595.4
$ws.Range("H24").Value = 2452.7058
$ws.Range("I24").Value = 1306
$ws.Range("J24").Value = 3255.4
$ws.Range("K24").Value = 1306
$ws.Range("L24").Value = 3255.4
$ws.Range("M24").Value = -1136
$ws.Range("N24").Value = -3595.4
$ws.Range("H70").Value = 29966
$ws.Range("J70").Value = 29966
$ws.Range("L70").Value = 29966
$ws.Range("N70").Value = -30596
$ws.Range("H73").Value = 29966
$ws.Range("J73").Value = 29966
$ws.Range("L73").Value = 29966
$ws.Range("N73").Value = -32150
$ws.Range("H113").Value = 4349.375
$ws.Range("I113").Value = 3874.25
$ws.Range("J113").Value = 4824.5
$ws.Range("K113").Value = 3874.25
$ws.Range("L113").Value = 4824.5
$ws.Range("M113").Value = -1704.25
$ws.Range("N113").Value = -9164.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 4998.5
$ws.Range("J19").Value = 4999
$ws.Range("L19").Value = 14997
$ws.Range("N19").Value = -15345
$ws.Range("H74").Value = 4495
$ws.Range("I74").Value = 4495
$ws.Range("K74").Value = 13485
$ws.Range("M74").Value = -12424
$ws.Range("H77").Value = 4495
$ws.Range("I77").Value = 4495
$ws.Range("K77").Value = 40455
$ws.Range("M77").Value = -35151
$ws.Range("H97").Value = 1852.6666
$ws.Range("J97").Value = 1973
$ws.Range("L97").Value = 5919
$ws.Range("N97").Value = -6911
$ws.Range("H117").Value = 1755.5
$ws.Range("I117").Value = 1135.75
$ws.Range("J117").Value = 2995
$ws.Range("K117").Value = 3407.25
$ws.Range("L117").Value = 8985
$ws.Range("M117").Value = 34.75
$ws.Range("H121").Value = 15152518
$ws.Range("I121").Value = 41666900
$ws.Range("J121").Value = 1442.4286
$ws.Range("K121").Value = 125000700
$ws.Range("L121").Value = 4327.2858
$ws.Range("M121").Value = -124999390
$ws.Range("N121").Value = -6947.2858
$ws.Range("H122").Value = 2026.1875
$ws.Range("J122").Value = 2298.111
$ws.Range("L122").Value = 20682.999
$ws.Range("N122").Value = -25582.999
$ws.Range("H133").Value = 1846
$ws.Range("I133").Value = 1846
$ws.Range("K133").Value = 5538
$ws.Range("M133").Value = -478
$ws.Range("H136").Value = 4079.3462
$ws.Range("I136").Value = 3582.6155
$ws.Range("K136").Value = 10747.8465
$ws.Range("M136").Value = -5647.8465
$ws.Range("H137").Value = 5963.409
$ws.Range("J137").Value = 10062.7
$ws.Range("L137").Value = 30188.1
$ws.Range("N137").Value = -40388.10000000001
$ws.Range("N117").Value = -15869

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 25000000
$ws.Range("J3").Value = 40000000
$ws.Range("L3").Value = 40000000
$ws.Range("N3").Value = -40000232

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1036.9231
$ws.Range("I107").Value = 1077.1666
$ws.Range("K107").Value = 3231.4998
$ws.Range("M107").Value = -1311.4998
$ws.Range("H132").Value = 1388.0605
$ws.Range("J132").Value = 2112.8
$ws.Range("L132").Value = 6338.400000000001
$ws.Range("N132").Value = -11398.4
